# supervision_en.xlsx -- add LaTeX \href{...}{...} hyperlink code around
# university / co-supervisor names, restyle the data range (wrap + top-left
# align), size the "what" rows to fit the now-longer wrapped text, extend the
# pre-formatted (but still empty) data range down to row 38, and refresh the
# sheet's saved scroll position / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rewrite the "where" (D) and "why" (E) column text to wrap the
#    university / collaborator names in vitae's \href{url}{text} macro.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = '\href{https://www.uv.es/}{Universitat de València}, España'
$ws.Range("E2").Value = '\href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita}. Supervised together with  Alicia Salvador'
$ws.Range("D3").Value = '\href{https://www.uel.ac.uk/}{University of East London}, UK'
$ws.Range("E3").Value = '\href{https://www.researchgate.net/profile/Francisco-Flores-14}{Francisco Javier Flores}. Supervised together with Lisa Chiara Fellin'
$ws.Range("D4").Value = '\href{https://www.stir.ac.uk/}{University of Stirling}, UK'
$ws.Range("E4").Value = 'Julia Sanz-Vidania. Supervised together with \href{https://www.scraigroberts.com/}{S Craig Roberts}'
$ws.Range("D5").Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia'
$ws.Range("E5").Value = 'Adrián Acosta Guerrero. Supervised together with \href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita}'
$ws.Range("D6").Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia'
$ws.Range("D7").Value = '\href{https://www.upn.edu.co/}{Universidad Pedagógica Nacional}, Colombia'
$ws.Range("D9").Value = '\href{https://www.unbosque.edu.co/}{Universidad El Bosque}, Colombia'

# ---------------------------------------------------------------------------
# 2. Extend the (still blank) data area down to row 38 so the sheet carries
#    pre-formatted rows ready for future entries.
# ---------------------------------------------------------------------------
$ws.Range("E27:E38").Value = ""

# ---------------------------------------------------------------------------
# 3. Re-style the whole when/with/where/why block: wrap text, left/top align.
#    Format a single already-wrapped cell first and fan it out with a style
#    copy so every cell in B1:E38 lands on one shared style.
# ---------------------------------------------------------------------------
$ws.Range("E2").HorizontalAlignment = -4131
$ws.Range("E2").VerticalAlignment = -4160
$ws.Range("E2").Copy() | Out-Null
$ws.Range("B1:E38").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Taller rows for the entries whose "where"/"why" text now wraps across
#    multiple lines.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45

# ---------------------------------------------------------------------------
# 5. Reset the saved view: scroll back to column A and leave the selection on
#    a single cell (E6) instead of the old A6:XFD8 block.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E6").Select() | Out-Null
